$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.56699999999999
$ws.Range("A6").Value = -22.70050000000002
$ws.Range("A7").Value = -21.995
$ws.Range("B7").Value = 5.1722
$ws.Range("A8").Value = -22.33570000000002
$ws.Range("B11").Value = 5.331400000000001
$ws.Range("B12").Value = 4.695499999999999
$ws.Range("D12").Value = -7.121899999999997
$ws.Range("D13").Value = -8.658799999999996
$ws.Range("D14").Value = -8.192599999999997
$ws.Range("B15").Value = 4.899799999999997
$ws.Range("A16").Value = -21.64369999999999
$ws.Range("D16").Value = -9.005100000000009
$ws.Range("D19").Value = -7.386499999999999
$ws.Range("A20").Value = -22.5488
$ws.Range("B20").Value = 4.366899999999996
$ws.Range("D20").Value = -7.681399999999994
$ws.Range("A21").Value = -22.2667
$ws.Range("B21").Value = 5.244199999999998
$ws.Range("B22").Value = 9.158600000000002
$ws.Range("D22").Value = -8.115900000000002
$ws.Range("B23").Value = 9.018900000000007
$ws.Range("A28").Value = -22.3389
$ws.Range("A29").Value = -21.7146
$ws.Range("B29").Value = 4.921800000000001
$ws.Range("A30").Value = -21.74630000000002
$ws.Range("A32").Value = -21.20689999999999
$ws.Range("B34").Value = 9.411800000000008
$ws.Range("D36").Value = -8.172999999999995
$ws.Range("A40").Value = -19.3115
$ws.Range("B42").Value = 9.319399999999996
$ws.Range("B43").Value = 5.115800000000003
$ws.Range("D43").Value = -8.724799999999997
$ws.Range("B44").Value = 5.584900000000001
$ws.Range("B45").Value = 4.954600000000002
$ws.Range("A46").Value = -22.19540000000001
$ws.Range("B46").Value = 5.627200000000001
$ws.Range("D46").Value = -8.038599999999995
$ws.Range("B50").Value = 4.392799999999997
$ws.Range("D50").Value = -8.199599999999997
$ws.Range("A51").Value = -22.242
$ws.Range("B51").Value = 5.221599999999999
$ws.Range("A52").Value = -21.95739999999999
$ws.Range("A57").Value = -22.61780000000001
$ws.Range("B57").Value = 5.347899999999997
$ws.Range("A59").Value = -22.1537
$ws.Range("A62").Value = -22.13120000000001
$ws.Range("B65").Value = 5.047900000000004
$ws.Range("A66").Value = -21.5521
$ws.Range("B66").Value = 5.012199999999998
$ws.Range("B67").Value = 5.290699999999999
$ws.Range("A73").Value = -20.03849999999998
$ws.Range("A74").Value = -21.96599999999998
$ws.Range("D76").Value = -8.065800000000001
$ws.Range("A77").Value = -20.03109999999998
$ws.Range("B79").Value = 9.332100000000004
$ws.Range("B84").Value = 5.534700000000002
$ws.Range("B87").Value = 5.011899999999996
$ws.Range("A92").Value = -21.46000000000002
$ws.Range("B92").Value = 4.832699999999995
$ws.Range("D95").Value = -8.097200000000001
$ws.Range("B97").Value = 6.214100000000003
$ws.Range("D97").Value = -8.473499999999994
$ws.Range("D99").Value = -8.2844
$ws.Range("A100").Value = -22.1388
